$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C, shifting T.C (Azure)/T.C (Desc.)/Error right
$ws.Columns("C:C").Insert()

# New header + value for the inserted PRINCIPAL column
$ws.Range("C1").Value = "PRINCIPAL"
$ws.Range("C2").Value = 10

# Match the new selection recorded in the sheet view
$ws.Range("C2").Select()
